{"js": "// Replace each arithmetic expression in the practice table with its new value.\n// Pairs are (old expression, new expression) taken from the commit diff. Each\n// old expression is unique in the document, so a body.search + insertText(replace)\n// swap keeps the surrounding run formatting (font, size) untouched.\nconst body = context.document.body;\n\nconst pairs = [\n  [\"43+28=\", \"81-47=\"],\n  [\"30+0=\", \"72-7=\"],\n  [\"28+69=\", \"43-25=\"],\n  [\"82-9=\", \"13+76=\"],\n  [\"5+20=\", \"19+5=\"],\n  [\"46+18=\", \"99-94=\"],\n  [\"24-6=\", \"47-45=\"],\n  [\"84-7=\", \"54+45=\"],\n  [\"55-3=\", \"74-23=\"],\n  [\"92-36=\", \"59-18=\"],\n  [\"3+65=\", \"27-23=\"],\n  [\"99-26=\", \"69+2=\"],\n  [\"26+24=\", \"9+25=\"],\n  [\"83-51=\", \"15+64=\"],\n  [\"99-91=\", \"5-4=\"],\n  [\"83-6=\", \"29-12=\"],\n  [\"85-55=\", \"74-37=\"],\n  [\"54-46=\", \"88-66=\"],\n  [\"94-70=\", \"2+26=\"],\n  [\"13-0=\", \"81-36=\"],\n  [\"3+18=\", \"50-18=\"],\n  [\"64+28=\", \"37-7=\"],\n  [\"70+27=\", \"77+7=\"],\n  [\"38+23=\", \"85+9=\"],\n  [\"67-2=\", \"73+3=\"],\n  [\"72-49=\", \"62+34=\"],\n  [\"37+47=\", \"68-3=\"],\n  [\"63-50=\", \"24+24=\"],\n  [\"11+68=\", \"45+38=\"],\n  [\"30-12=\", \"41-29=\"],\n  [\"53-30=\", \"33-8=\"],\n  [\"56+16=\", \"39+4=\"],\n  [\"63-38=\", \"66-18=\"],\n  [\"42+38=\", \"27+58=\"],\n  [\"12+82=\", \"68-12=\"],\n  [\"94-36=\", \"23+41=\"],\n  [\"95-65=\", \"34+30=\"],\n  [\"14+77=\", \"63-16=\"],\n  [\"80-36=\", \"56+18=\"],\n  [\"4+83=\", \"70+12=\"],\n  [\"11+75=\", \"34-27=\"],\n  [\"30+53=\", \"43-4=\"],\n  [\"69-30=\", \"77-38=\"],\n  [\"36+20=\", \"27+28=\"],\n  [\"44-35=\", \"49+33=\"],\n  [\"5+86=\", \"56-3=\"],\n  [\"16+11=\", \"41+26=\"],\n  [\"43+45=\", \"89-29=\"],\n  [\"35-18=\", \"39-34=\"],\n  [\"86-77=\", \"11+13=\"],\n  [\"23+31=\", \"51-50=\"],\n  [\"42-9=\", \"44+25=\"],\n  [\"7+30=\", \"96-15=\"],\n  [\"87-64=\", \"29+36=\"],\n  [\"3+60=\", \"43+35=\"],\n  [\"56+34=\", \"44-2=\"],\n  [\"84-21=\", \"65-27=\"],\n  [\"26+15=\", \"74-70=\"],\n  [\"51-25=\", \"54-22=\"],\n  [\"58+28=\", \"4+14=\"],\n  [\"14+7=\", \"96-10=\"],\n  [\"11+4=\", \"67-16=\"],\n  [\"20+32=\", \"45+35=\"],\n  [\"56-51=\", \"15+70=\"],\n  [\"60-39=\", \"79-27=\"],\n  [\"69-54=\", \"56-15=\"],\n  [\"90-66=\", \"2+62=\"],\n  [\"11+70=\", \"49-30=\"],\n  [\"81-26=\", \"99-86=\"],\n  [\"80-5=\", \"1+28=\"],\n  [\"40+59=\", \"22+65=\"],\n  [\"55+16=\", \"8+1=\"],\n  [\"52-48=\", \"52-26=\"],\n  [\"29-4=\", \"15+60=\"],\n  [\"91-44=\", \"52+31=\"],\n  [\"19+4=\", \"15+59=\"],\n  [\"0+8=\", \"56-39=\"],\n  [\"46+24=\", \"0+3=\"],\n  [\"79+11=\", \"83+3=\"],\n  [\"21+37=\", \"66+12=\"],\n  [\"95-52=\", \"52+3=\"],\n  [\"43+34=\", \"66-13=\"],\n  [\"67-37=\", \"22+0=\"],\n  [\"85-81=\", \"17+47=\"],\n  [\"41+21=\", \"54-53=\"],\n  [\"3+27=\", \"97-39=\"],\n  [\"23+34=\", \"93-5=\"],\n  [\"76-47=\", \"94-8=\"],\n  [\"30-13=\", \"85-68=\"],\n  [\"85-27=\", \"86-7=\"],\n  [\"19+79=\", \"16+5=\"],\n  [\"6+10=\", \"38+51=\"],\n  [\"66+15=\", \"89-1=\"],\n  [\"19+32=\", \"64-52=\"],\n  [\"74+7=\", \"24+74=\"],\n  [\"68-38=\", \"3+42=\"],\n  [\"37+33=\", \"99-42=\"],\n  [\"47-31=\", \"67-0=\"],\n  [\"13+85=\", \"87+3=\"],\n  [\"81-5=\", \"83-10=\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each arithmetic expression in the practice table with its new value.\n# Pairs are (old expression, new expression) taken from the commit diff, applied\n# in document order via Find/Replace so run formatting (font, size) is preserved.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"43+28=\", \"81-47=\")\n    ,@(\"30+0=\", \"72-7=\")\n    ,@(\"28+69=\", \"43-25=\")\n    ,@(\"82-9=\", \"13+76=\")\n    ,@(\"5+20=\", \"19+5=\")\n    ,@(\"46+18=\", \"99-94=\")\n    ,@(\"24-6=\", \"47-45=\")\n    ,@(\"84-7=\", \"54+45=\")\n    ,@(\"55-3=\", \"74-23=\")\n    ,@(\"92-36=\", \"59-18=\")\n    ,@(\"3+65=\", \"27-23=\")\n    ,@(\"99-26=\", \"69+2=\")\n    ,@(\"26+24=\", \"9+25=\")\n    ,@(\"83-51=\", \"15+64=\")\n    ,@(\"99-91=\", \"5-4=\")\n    ,@(\"83-6=\", \"29-12=\")\n    ,@(\"85-55=\", \"74-37=\")\n    ,@(\"54-46=\", \"88-66=\")\n    ,@(\"94-70=\", \"2+26=\")\n    ,@(\"13-0=\", \"81-36=\")\n    ,@(\"3+18=\", \"50-18=\")\n    ,@(\"64+28=\", \"37-7=\")\n    ,@(\"70+27=\", \"77+7=\")\n    ,@(\"38+23=\", \"85+9=\")\n    ,@(\"67-2=\", \"73+3=\")\n    ,@(\"72-49=\", \"62+34=\")\n    ,@(\"37+47=\", \"68-3=\")\n    ,@(\"63-50=\", \"24+24=\")\n    ,@(\"11+68=\", \"45+38=\")\n    ,@(\"30-12=\", \"41-29=\")\n    ,@(\"53-30=\", \"33-8=\")\n    ,@(\"56+16=\", \"39+4=\")\n    ,@(\"63-38=\", \"66-18=\")\n    ,@(\"42+38=\", \"27+58=\")\n    ,@(\"12+82=\", \"68-12=\")\n    ,@(\"94-36=\", \"23+41=\")\n    ,@(\"95-65=\", \"34+30=\")\n    ,@(\"14+77=\", \"63-16=\")\n    ,@(\"80-36=\", \"56+18=\")\n    ,@(\"4+83=\", \"70+12=\")\n    ,@(\"11+75=\", \"34-27=\")\n    ,@(\"30+53=\", \"43-4=\")\n    ,@(\"69-30=\", \"77-38=\")\n    ,@(\"36+20=\", \"27+28=\")\n    ,@(\"44-35=\", \"49+33=\")\n    ,@(\"5+86=\", \"56-3=\")\n    ,@(\"16+11=\", \"41+26=\")\n    ,@(\"43+45=\", \"89-29=\")\n    ,@(\"35-18=\", \"39-34=\")\n    ,@(\"86-77=\", \"11+13=\")\n    ,@(\"23+31=\", \"51-50=\")\n    ,@(\"42-9=\", \"44+25=\")\n    ,@(\"7+30=\", \"96-15=\")\n    ,@(\"87-64=\", \"29+36=\")\n    ,@(\"3+60=\", \"43+35=\")\n    ,@(\"56+34=\", \"44-2=\")\n    ,@(\"84-21=\", \"65-27=\")\n    ,@(\"26+15=\", \"74-70=\")\n    ,@(\"51-25=\", \"54-22=\")\n    ,@(\"58+28=\", \"4+14=\")\n    ,@(\"14+7=\", \"96-10=\")\n    ,@(\"11+4=\", \"67-16=\")\n    ,@(\"20+32=\", \"45+35=\")\n    ,@(\"56-51=\", \"15+70=\")\n    ,@(\"60-39=\", \"79-27=\")\n    ,@(\"69-54=\", \"56-15=\")\n    ,@(\"90-66=\", \"2+62=\")\n    ,@(\"11+70=\", \"49-30=\")\n    ,@(\"81-26=\", \"99-86=\")\n    ,@(\"80-5=\", \"1+28=\")\n    ,@(\"40+59=\", \"22+65=\")\n    ,@(\"55+16=\", \"8+1=\")\n    ,@(\"52-48=\", \"52-26=\")\n    ,@(\"29-4=\", \"15+60=\")\n    ,@(\"91-44=\", \"52+31=\")\n    ,@(\"19+4=\", \"15+59=\")\n    ,@(\"0+8=\", \"56-39=\")\n    ,@(\"46+24=\", \"0+3=\")\n    ,@(\"79+11=\", \"83+3=\")\n    ,@(\"21+37=\", \"66+12=\")\n    ,@(\"95-52=\", \"52+3=\")\n    ,@(\"43+34=\", \"66-13=\")\n    ,@(\"67-37=\", \"22+0=\")\n    ,@(\"85-81=\", \"17+47=\")\n    ,@(\"41+21=\", \"54-53=\")\n    ,@(\"3+27=\", \"97-39=\")\n    ,@(\"23+34=\", \"93-5=\")\n    ,@(\"76-47=\", \"94-8=\")\n    ,@(\"30-13=\", \"85-68=\")\n    ,@(\"85-27=\", \"86-7=\")\n    ,@(\"19+79=\", \"16+5=\")\n    ,@(\"6+10=\", \"38+51=\")\n    ,@(\"66+15=\", \"89-1=\")\n    ,@(\"19+32=\", \"64-52=\")\n    ,@(\"74+7=\", \"24+74=\")\n    ,@(\"68-38=\", \"3+42=\")\n    ,@(\"37+33=\", \"99-42=\")\n    ,@(\"47-31=\", \"67-0=\")\n    ,@(\"13+85=\", \"87+3=\")\n    ,@(\"81-5=\", \"83-10=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
